$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# A tiny epsilon (in points) compensates for the single-precision float
# round-trip PowerPoint's Left/Top/Width/Height setters do internally
# before re-quantizing back to EMUs, so the stored EMU values land on
# the exact target instead of being truncated one EMU short.
$eps = 0.00002

# "TextBox 33" - [More Games] -> [Series Over]
$s.Shapes.Item(19).TextFrame.TextRange.Text = "[Series Over]"

# "TextBox 37" - reposition/resize, and split "Play Match" into two
# centered paragraphs: "Play Next" / "Match"
$playNext = $s.Shapes.Item(22)
$playNext.Left = (8813923 / 12700) + $eps
$playNext.Top = (2958827 / 12700) + $eps
$playNext.Width = (1687679 / 12700) + $eps
$playNext.Height = (646331 / 12700) + $eps
$playNext.TextFrame.TextRange.Text = "Play Next" + [char]13 + "Match"
$playNext.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# "TextBox 41" - "Yes" label, reposition only
$yes = $s.Shapes.Item(23)
$yes.Left = (7124264 / 12700) + $eps
$yes.Top = (3471805 / 12700) + $eps

# "TextBox 43" - "No" label, reposition only
$no = $s.Shapes.Item(24)
$no.Left = (8327814 / 12700) + $eps
$no.Top = (2912660 / 12700) + $eps
